$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell values
$ws.Range("B2").Value = 114
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 50

# Remove rows 4 and 5 entirely
$ws.Rows("4:5").Delete()
